# Applies the "Updated metadata and technotes" commit to the Technical
# Guidance workbook:
#  - Rewords several notes/technote cells on the "Quantity" sheet to
#    reflect pupil-number comparisons as "change" (rather than
#    "increase"/"growth") between 2009/10, 2022/23 and 2024/25, and to be
#    explicit about the January census dates used.
#  - Moves the active-cell selection on the "Quantity" sheet from B13 to
#    B12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quantity")

# Row 12 - "Change in pupil numbers" indicator block.
$ws.Range("A12").Value = "Change in pupil numbers`n2009/10 to 2022/23; Anticipated change in pupil numbers 2022/23 to 2024/25"

$ws.Range("B12").Value = "The actual percentage change in pupil numbers between 2009/10 and 2022/23; the anticipated percentage change in pupil numbers in primary or secondary state-funded mainstream provision between the 2022/23 and 2024/25 academic years."

$ws.Range("E12").Value = "1. Number of pupils in roll in January 2010, including dual registrations, in reception to year 11 in the following types of schools: Academy converter, Academy sponsor led, City technology college. Community school, Foundation school, Free schools, Studio schools, University technical college, Voluntary aided school, Voluntary controlled school. "

# Row 13 - second half of the same indicator block (2022/23 pupil numbers).
$ws.Range("C13").Value = "2. Pupil Numbers for the 2022/23 academic year taken from the pupil census in January 2023"

$ws.Range("E13").Value = "1. Number of pupils in roll in January 2023, including dual registrations, in reception to year 11 in the following types of schools: Academy converter, Academy sponsor led, City technology college. Community school, Foundation school, Free schools, Studio schools, University technical college, Voluntary aided school, Voluntary controlled school. `t`t`t`t`t`t`t"

# Restore the author's active selection (was B13, now B12) on the
# "Quantity" sheet.
$ws.Activate()
$ws.Range("B12").Select()
